# Apply data updates to rows 2-11 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 190;   B = "Gabrielly Ribeiro";        C = "Operações";              D = "Problemas pessoais"; E = 4; F = 45105; G = 9615.809999999999 },
    @{ Row = 3;  A = 98933; B = "Sr. Vitor Hugo Gonçalves";  C = "Atendimento ao Cliente";  D = "Viagem de negócios"; E = 1; F = 45101; G = 3737.68 },
    @{ Row = 4;  A = 7300;  B = "Evelyn Barros";             C = "Atendimento ao Cliente";  D = "Consulta médica";    E = 3; F = 45093; G = 3442.63 },
    @{ Row = 5;  A = 3750;  B = "Enzo Gabriel Dias";         C = "Jurídico";                D = "Doença";             E = 6; F = 45086; G = 3956.8 },
    @{ Row = 6;  A = 61629; B = "Pietro Oliveira";           C = "Engenharia";              D = "Outros";             E = 5; F = 45084; G = 12346.19 },
    @{ Row = 7;  A = 38013; B = "Lucas Gabriel da Costa";    C = "Vendas";                  D = "Problemas pessoais"; E = 7; F = 45102; G = 6444.36 },
    @{ Row = 8;  A = 93248; B = "Anthony da Mota";           C = "P&D";                     D = "Viagem de negócios"; E = 7; F = 45096; G = 3717.89 },
    @{ Row = 9;  A = 73383; B = "Ana Sophia Lima";           C = "Recursos Humanos";        D = "Problemas pessoais"; E = 7; F = 45103; G = 6358.84 },
    @{ Row = 10; A = 11548; B = "Pedro Miguel Alves";        C = "Operações";               D = "Problemas pessoais"; E = 1; F = 45085; G = 8803.75 },
    @{ Row = 11; A = 9181;  B = "Vitória Fernandes";         C = "Marketing";               D = "Outros";             E = 8; F = 45091; G = 3246.43 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
